# auto:removing today's date from viral and referral form
#
# The "date1" question (type=date, name=date1, label="Today's Date") on the
# "survey" sheet stops auto-calculating/locking itself to today's date: its
# `calculation` (today()) and `readonly` (TRUE) cells are removed so it is
# no longer a hidden, auto-filled field.
#
# The "select_one follow" question (name=patient, label="Patient was
# referred for a visit") is removed entirely from the survey.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# Drop the calculation/readonly cells for the "date1" row (row 28) so the
# field becomes a normal (non-auto-calculated) date prompt.
$ws.Range("F28:G28").Clear()

# Remove the whole "select_one follow" row (row 29) -- patient referral
# confirmation question -- shifting the rows below it up.
$ws.Range("A29:C29").EntireRow.Delete()
